$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$data = @(
    ,("lab.vape.table.setup", "Setup")
    ,("lab.vape.table.mixture", "Mix")
    ,("lab.vape.table.rating", "Hodnocení")
    ,("lab.vape.table.taste", "Chuť")
    ,("lab.vape.table.power", "Výkon")
    ,("lab.vape.table.tc", "Teplota")
    ,("lab.vape.preview.setup", "Setup")
    ,("lab.vape.preview.mixture", "Mix")
    ,("lab.vape.preview.driptip", "Náústek")
    ,("lab.vape.preview.leaks", "Úniky")
    ,("lab.vape.preview.dryhit", "Dryhity")
    ,("lab.vape.preview.rating", "Celkové hodnocení")
    ,("lab.vape.preview.taste", "Hodnocení chuti")
    ,("lab.vape.preview.power", "Výkon (watty)")
    ,("lab.vape.preview.tc", "Teplota")
    ,("lab.vape.preview.airflow", "Airflow")
    ,("lab.vape.preview.juice", "Juice flow")
    ,("lab.vape.preview.mtl", "MTL")
    ,("lab.vape.preview.dl", "DL")
    ,("lab.vape.preview.clouds", "Oblaka")
    ,("lab.vape.preview.fruits", "Ovocné tóny")
    ,("lab.vape.preview.tobacco", "Tabák")
    ,("lab.vape.preview.cakes", "Buchty")
    ,("lab.vape.preview.complex", "Komplexní")
    ,("lab.vape.preview.fresh", "Větrnost")
    ,("lab.vape.preview.atomizer", "Atomizér")
    ,("lab.vape.preview.coil", "Spirálka")
    ,("lab.vape.preview.mod", "Mod")
)

$startRow = 464
$templateRow = $ws.Range("A463:C463")
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $destRow = $ws.Range("A" + $row + ":C" + $row)
    $templateRow.Copy()
    $destRow.PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = "cs"
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

$ws.Application.GoTo($ws.Range("A478"), $true) | Out-Null
$ws.Range("B485").Select() | Out-Null
